$d = $word.ActiveDocument

function Split-At($startPos, $pairs) {
    # $pairs is an array of @(relStart, relEnd) pairs (relative to $startPos);
    # toggling Bold on/off on each one is a no-visual-effect way to force a
    # run boundary at those exact offsets.
    foreach ($pair in $pairs) {
        $a = $startPos + $pair[0]
        $b = $startPos + $pair[1]
        if ($b -gt $a) {
            $seg = $d.Range($a, $b)
            $seg.Font.Bold = 1
            $seg.Font.Bold = 0
        }
    }
}

# ===========================================================================
# Text edits (content-level) -- these are the actual wording changes.
# ===========================================================================

# Edit 1: "show a list of" -> "display"
$eA = $d.Content
$eA.Find.Execute("show a list of", $true, $false, $false, $false, $false, `
                  $true, 1, $false, "display", 2)

# Edit 3: ", show the list of countries with the highest coefficient." ->
#         ", display the countries with the highest coefficient."
$eB = $d.Content
$eB.Find.Execute(", show the list of countries with the highest coefficient.", `
                  $true, $false, $false, $false, $false, $true, 1, $false, `
                  ", display the countries with the highest coefficient.", 2)

# Edit 4: merge "...cumulative scores c" + bookmark + "ompare?" into one run
#         with the bookmark removed from this location.
$eC = $d.Content
$eC.Find.Execute("5 – How do the USSR and Russia’s cumulative scores compare?", `
                  $true, $false, $false, $false, $false, $true, 1, $false, `
                  "5 – How do the USSR and Russia’s cumulative scores compare?", 2)

# ===========================================================================
# Paragraph A: the "1 – Browse ..." numbered block (contains edits 1 and 3,
# and also gains the relocated bookmark around "using the derivat|ive").
# ===========================================================================
$findA = $d.Content
$findA.Find.Execute("Browse", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$paraA = $findA.Paragraphs(1)
$startA = $paraA.Range.Start

$pairsA = @(
    @(0,1), @(1,4), @(4,11), @(11,12), @(12,13), @(13,20), @(20,74), @(74,89), @(89,90),
    @(91,104), @(104,105), @(105,168),
    @(169,170), @(170,173), @(173,180), @(180,181), @(181,240),
    @(241,242), @(242,245), @(245,255), @(255,272), @(272,285), @(285,297), @(297,298),
    @(298,310), @(310,311), @(311,319), @(319,320), @(320,324), @(324,363),
    @(364,365), @(365,368), @(368,387), @(387,411)
)
Split-At $startA $pairsA

# Relocate the "_GoBack" bookmark to sit between "using the derivat" and "ive variable "
$oldBm = $d.Bookmarks.Item("_GoBack")
$oldBm.Delete()
$bmPos = $startA + 272
$bmRng = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRng)

# ===========================================================================
# Paragraph B: the "1 – What countries..." numbered block (ends with the now
# merged "5 – How do the USSR and Russia’s cumulative scores compare?").
# ===========================================================================
$findB = $d.Content
$findB.Find.Execute("What countries had the most gold", $true, $false, $false, $false, $false, `
                     $true, 1, $false, "", 0)
$paraB = $findB.Paragraphs(1)
$startB = $paraB.Range.Start

$pairsB = @(
    @(0,1), @(1,4), @(4,37), @(37,68), @(68,76),
    @(77,78), @(78,81), @(81,126),
    @(127,128), @(128,131), @(131,174),
    @(175,176), @(176,179), @(179,196), @(196,201), @(201,241),
    @(242,301)
)
Split-At $startB $pairsB
